$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.629.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.243.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.21%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.64%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -1.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.84%  "

$ws.Range("E11").Value = "  -1.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.43%  "

$ws.Range("E13").Value = "  +0.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.586.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.242.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.831"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.392.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0936"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.70%  "

$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.21%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("E24").Value = "  -4.66%  "

$ws.Range("E25").Value = "  -1.19%  "

$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("E27").Value = "  +4.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0782"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.34%  "

$ws.Range("E34").Value = "  +0.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.109"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.27%  "

$ws.Range("E37").Value = "  -1.56%  "

$ws.Range("E38").Value = "  +5.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.77%  "

$ws.Range("E42").Value = "  +0.48%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.808.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "82.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.48%  "

$ws.Range("E47").Value = "  -1.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "97.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.52%  "

$ws.Range("E49").Value = "  -1.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.54%  "
